# Update mods data [2026-01-06 15:11:46]
# Append the next day's mod-count reading as a new row at the bottom of
# the ModCounts sheet (row 57): Date, Game, ModCount.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 57

$date = "2026/01/06"
$game = "逃离鸭科夫"
$modCount = 1138

# Column A holds the date as literal text (matching every other row in
# the sheet, which stores "yyyy/mm/dd" as a string rather than a real
# date serial). Assigning the string straight to .Value would make Excel
# "smart"-parse it into a date serial, so instead we build it via a
# formula (never subject to that auto-detection) and then flatten the
# formula down to a plain value with a copy / paste-values round-trip.
$cellA = $ws.Cells.Item($newRow, 1)
$cellA.Formula = '="' + $date + '"'
$cellA.Copy()
$cellA.PasteSpecial(-4163)  # xlPasteValues

$cellB = $ws.Cells.Item($newRow, 2)
$cellB.Value = $game

$cellC = $ws.Cells.Item($newRow, 3)
$cellC.Value = $modCount

# Match the centered alignment style ("s=1") used by every other data
# row in the sheet.
foreach ($col in 1..3) {
    $cell = $ws.Cells.Item($newRow, $col)
    $cell.HorizontalAlignment = -4108  # xlCenter
    $cell.VerticalAlignment = -4108    # xlCenter
}
